# Update the "Training Dashboard" sheet with the new progress date (04-Nov-2025).
# For each data row (3-9):
#   - Column H (PERIOD TO EXPIRE) decreases by 1 day
#   - Column I (LAST UPDATE) changes from "03-Nov-2025" to "04-Nov-2025"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$rows = 3..9

foreach ($r in $rows) {
    # PERIOD TO EXPIRE (col H) drops by one day.
    $hCell = $ws.Cells.Item($r, 8)
    $current = $hCell.Value2
    $hCell.Value = $current - 1

    # LAST UPDATE (col I) moves from 03-Nov-2025 to 04-Nov-2025.
    # Force the cell to stay plain text (matches source data) instead of
    # letting Excel auto-convert the recognizable date string into a date
    # serial number.
    $iCell = $ws.Cells.Item($r, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
}
